$wb = $excel.ActiveWorkbook

# --- Update the "Date" property on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-23T20:31:36+00:00"

# --- Update the "Code" column (B) on the Concepts sheet so it matches ---
# --- the "Display" column (C) for each concept row (the hyphenated    ---
# --- code strings are replaced by their spaced display equivalents).  ---
$concepts = $wb.Worksheets.Item("Concepts")

for ($row = 2; $row -le 9; $row++) {
    $display = $concepts.Cells.Item($row, 3).Value2
    $concepts.Cells.Item($row, 2).Value = $display
}
